$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("tarea1")
$ws.Range("B2:B27").NumberFormat = "@"
$ws.Range("B2").Value = "9.81"
$ws.Range("B3").Value = "16.09"
$ws.Range("B4").Value = "22.37"
$ws.Range("B5").Value = "28.65"
$ws.Range("B6").Value = "34.93"
$ws.Range("B7").Value = "41.21"
$ws.Range("B8").Value = "47.49"
$ws.Range("B9").Value = "53.77"
$ws.Range("B10").Value = "60.05"
$ws.Range("B11").Value = "66.33"
$ws.Range("B12").Value = "72.61"
$ws.Range("B13").Value = "78.89"
$ws.Range("B14").Value = "85.17"
$ws.Range("B15").Value = "91.45"
$ws.Range("B16").Value = "97.73"
$ws.Range("B17").Value = "104.01"
$ws.Range("B18").Value = "110.29"
$ws.Range("B19").Value = "116.57"
$ws.Range("B20").Value = "122.85"
$ws.Range("B21").Value = "129.13"
$ws.Range("B22").Value = "135.41"
$ws.Range("B23").Value = "141.69"
$ws.Range("B24").Value = "147.97"
$ws.Range("B25").Value = "154.25"
$ws.Range("B26").Value = "160.53"
$ws.Range("B27").Value = "166.81"

$ws = $wb.Worksheets.Item("tarea2")
$ws.Range("B2:B27").NumberFormat = "@"
$ws.Range("B2").Value = "34.64"
$ws.Range("B3").Value = "25.59"
$ws.Range("B4").Value = "28.52"
$ws.Range("B5").Value = "40.42"
$ws.Range("B6").Value = "11.12"
$ws.Range("B7").Value = "61.20"
$ws.Range("B8").Value = "44.10"
$ws.Range("B9").Value = "73.33"
$ws.Range("B10").Value = "56.14"
$ws.Range("B11").Value = "57.83"
$ws.Range("B12").Value = "50.74"
$ws.Range("B13").Value = "90.15"
$ws.Range("B14").Value = "100.63"
$ws.Range("B15").Value = "83.65"
$ws.Range("B16").Value = "103.47"
$ws.Range("B17").Value = "89.48"
$ws.Range("B18").Value = "128.12"
$ws.Range("B19").Value = "140.35"
$ws.Range("B20").Value = "137.97"
$ws.Range("B21").Value = "122.06"
$ws.Range("B22").Value = "118.24"
$ws.Range("B23").Value = "146.39"
$ws.Range("B24").Value = "138.15"
$ws.Range("B25").Value = "158.85"
$ws.Range("B26").Value = "135.65"
$ws.Range("B27").Value = "187.84"

$ws = $wb.Worksheets.Item("tarea3")
$ws.Range("B2:B27").NumberFormat = "@"
$ws.Range("B2").Value = "-3.66"
$ws.Range("B3").Value = "93.89"
$ws.Range("B4").Value = "3.04"
$ws.Range("B5").Value = "202.57"
$ws.Range("B6").Value = "217.01"
$ws.Range("B7").Value = "287.28"
$ws.Range("B8").Value = "378.04"
$ws.Range("B9").Value = "693.23"
$ws.Range("B10").Value = "852.52"
$ws.Range("B11").Value = "1005.74"
$ws.Range("B12").Value = "1202.23"
$ws.Range("B13").Value = "1559.58"
$ws.Range("B14").Value = "1911.75"
$ws.Range("B15").Value = "2142.29"
$ws.Range("B16").Value = "2519.12"
$ws.Range("B17").Value = "2799.50"
$ws.Range("B18").Value = "3148.96"
$ws.Range("B19").Value = "3650.77"
$ws.Range("B20").Value = "4003.08"
$ws.Range("B21").Value = "4628.97"
$ws.Range("B22").Value = "5110.69"
$ws.Range("B23").Value = "5482.02"
$ws.Range("B24").Value = "6187.39"
$ws.Range("B25").Value = "6682.75"
$ws.Range("B26").Value = "7270.90"
$ws.Range("B27").Value = "7829.40"

$ws = $wb.Worksheets.Item("tarea4")
$ws.Range("B2:B27").NumberFormat = "@"
$ws.Range("B2").Value = "67.85"
$ws.Range("B3").Value = "35.37"
$ws.Range("B4").Value = "27.60"
$ws.Range("B5").Value = "32.37"
$ws.Range("B6").Value = "95.57"
$ws.Range("B7").Value = "65.69"
$ws.Range("B8").Value = "53.22"
$ws.Range("B9").Value = "72.80"
$ws.Range("B10").Value = "6.83"
$ws.Range("B11").Value = "39.99"
$ws.Range("B12").Value = "7.94"
$ws.Range("B13").Value = "92.95"
$ws.Range("B14").Value = "15.19"
$ws.Range("B15").Value = "25.39"
$ws.Range("B16").Value = "5.38"
$ws.Range("B17").Value = "10.83"
$ws.Range("B18").Value = "38.88"
$ws.Range("B19").Value = "77.00"
$ws.Range("B20").Value = "56.19"
$ws.Range("B21").Value = "37.41"
$ws.Range("B22").Value = "25.80"
$ws.Range("B23").Value = "90.98"
$ws.Range("B24").Value = "50.51"
$ws.Range("B25").Value = "18.94"
$ws.Range("B26").Value = "34.17"
$ws.Range("B27").Value = "39.74"
